# Apply scheduled-runner market data updates to the Bahamut_Profits workbook.
# Each block targets one leve row (identified by sheet + item name) and updates
# the currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) to the new
# values captured by the data-refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 17: One for the Road
$ws.Range("H17").Value = 237212.19
$ws.Range("J17").Value = 237212.19
$ws.Range("L17").Value = 711636.5700000001
$ws.Range("N17").Value = -711972.5700000001

# ALC row 33: Glazed and Confused
$ws.Range("H33").Value = 62593.5
$ws.Range("I33").Value = 66759.734
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 66759.734
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = -66530.734
$ws.Range("N33").Value = -558

# ALC row 74: Adhesive of Antipathy
$ws.Range("H74").Value = 4632263
$ws.Range("I74").Value = 5660432.5
$ws.Range("J74").Value = 5500
$ws.Range("K74").Value = 5660432.5
$ws.Range("L74").Value = 5500
$ws.Range("M74").Value = -5659496.5
$ws.Range("N74").Value = -7372

# ALC row 76: Warding Off Temptation
$ws.Range("H76").Value = 28951258
$ws.Range("I76").Value = 37934756
$ws.Range("J76").Value = 4432
$ws.Range("K76").Value = 37934756
$ws.Range("L76").Value = 4432
$ws.Range("M76").Value = -37934441
$ws.Range("N76").Value = -5062

# ALC row 77: It's Gonna Grow Back (L)
$ws.Range("H77").Value = 4632263
$ws.Range("I77").Value = 5660432.5
$ws.Range("J77").Value = 5500
$ws.Range("K77").Value = 28302162.5
$ws.Range("L77").Value = 27500
$ws.Range("M77").Value = -28297482.5
$ws.Range("N77").Value = -36860

# ALC row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 28951258
$ws.Range("I79").Value = 37934756
$ws.Range("J79").Value = 4432
$ws.Range("K79").Value = 37934756
$ws.Range("L79").Value = 4432
$ws.Range("M79").Value = -37933664
$ws.Range("N79").Value = -6616

# ALC row 92: Whinier than the Sword
$ws.Range("H92").Value = 1092.3704
$ws.Range("I92").Value = 1268.5
$ws.Range("J92").Value = 589.1429000000001
$ws.Range("K92").Value = 1268.5
$ws.Range("L92").Value = 589.1429000000001
$ws.Range("M92").Value = -20.5
$ws.Range("N92").Value = -3085.1429

# ALC row 96: Scroll Down
$ws.Range("H96").Value = 5019.1816
$ws.Range("I96").Value = 5213.9375
$ws.Range("J96").Value = 4499.8335
$ws.Range("K96").Value = 15641.8125
$ws.Range("L96").Value = 13499.5005
$ws.Range("M96").Value = -14268.8125
$ws.Range("N96").Value = -16245.5005

# ALC row 129: Practical Command
$ws.Range("H129").Value = 1684439.4
$ws.Range("J129").Value = 2058692.5
$ws.Range("L129").Value = 6176077.5
$ws.Range("N129").Value = -6186077.5

# ALC row 136: I Like Big Brush and I Cannot Lie
$ws.Range("H136").Value = 28850
$ws.Range("J136").Value = 28850
$ws.Range("L136").Value = 28850
$ws.Range("N136").Value = -39050

# ALC row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1636.3636
$ws.Range("J137").Value = 1666.6666
$ws.Range("L137").Value = 4999.9998
$ws.Range("N137").Value = -10099.9998

# ALC row 139: Something Salty and Ceremonial
$ws.Range("H139").Value = 34212
$ws.Range("J139").Value = 34212
$ws.Range("L139").Value = 34212
$ws.Range("N139").Value = -44492

# ALC row 140: Tome for Tradition
$ws.Range("H140").Value = 48020
$ws.Range("J140").Value = 48020
$ws.Range("L140").Value = 48020
$ws.Range("N140").Value = -58380

$ws = $wb.Worksheets.Item("ARM")
# ARM row 4: Eyes Bigger than the Plate
$ws.Range("H4").Value = 149
$ws.Range("I4").Value = 149
$ws.Range("K4").Value = 149
$ws.Range("M4").Value = -33

# ARM row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 3520.889
$ws.Range("I61").Value = 3520.889
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3520.889
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3308.889
$ws.Range("N61").ClearContents()

# ARM row 74: As the Bolt Flies
$ws.Range("H74").Value = 864.73334
$ws.Range("I74").Value = 1104.2
$ws.Range("J74").Value = 745
$ws.Range("K74").Value = 1104.2
$ws.Range("L74").Value = 745
$ws.Range("M74").Value = -230.2
$ws.Range("N74").Value = -2493

# ARM row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 864.73334
$ws.Range("I77").Value = 1104.2
$ws.Range("J77").Value = 745
$ws.Range("K77").Value = 5521
$ws.Range("L77").Value = 3725
$ws.Range("M77").Value = -1153
$ws.Range("N77").Value = -12461

# ARM row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1888.275
$ws.Range("I132").Value = 1462.7941
$ws.Range("J132").Value = 4299.3335
$ws.Range("K132").Value = 4388.3823
$ws.Range("L132").Value = 12898.0005
$ws.Range("M132").Value = -1858.3823
$ws.Range("N132").Value = -17958.0005

# ARM row 136: Metal with Mettle
$ws.Range("H136").Value = 3520.889
$ws.Range("I136").Value = 3520.889
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10562.667
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8012.667000000001
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# BSM row 132: Always Be Prepaired
$ws.Range("H132").Value = 100285.71
$ws.Range("J132").Value = 100285.71
$ws.Range("L132").Value = 100285.71
$ws.Range("N132").Value = -110405.71

# BSM row 134: Ruthenium Supremium
$ws.Range("H134").Value = 137253
$ws.Range("I134").Value = 4708.636
$ws.Range("J134").Value = 501750
$ws.Range("K134").Value = 14125.908
$ws.Range("L134").Value = 1505250
$ws.Range("M134").Value = -11590.908
$ws.Range("N134").Value = -1510320

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31: Wall Not Found
$ws.Range("H31").Value = 2632.6875
$ws.Range("I31").Value = 1641.7693
$ws.Range("J31").Value = 6926.6665
$ws.Range("K31").Value = 1641.7693
$ws.Range("L31").Value = 6926.6665
$ws.Range("M31").Value = -1346.7693
$ws.Range("N31").Value = -7516.6665

# CRP row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2632.6875
$ws.Range("I34").Value = 1641.7693
$ws.Range("J34").Value = 6926.6665
$ws.Range("K34").Value = 1641.7693
$ws.Range("L34").Value = 6926.6665
$ws.Range("M34").Value = -1439.7693
$ws.Range("N34").Value = -7330.6665

# CRP row 47: Grippy When Wet
$ws.Range("H47").Value = 11000
$ws.Range("I47").Value = 11000
$ws.Range("K47").Value = 11000
$ws.Range("M47").Value = -10434

# CRP row 140: Spear Pressure
$ws.Range("H140").Value = 54446.668
$ws.Range("J140").Value = 54446.668
$ws.Range("L140").Value = 54446.668
$ws.Range("N140").Value = -64806.668

$ws = $wb.Worksheets.Item("CUL")
# CUL row 46: Feeding Frenzy
$ws.Range("H46").Value = 1999
$ws.Range("J46").Value = 1999
$ws.Range("L46").Value = 5997
$ws.Range("N46").Value = -6179

# CUL row 113: Can't Eat Just One
$ws.Range("H113").Value = 584.13794
$ws.Range("I113").Value = 550
$ws.Range("J113").Value = 586.6667
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 1760.0001
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -6100.0001

# CUL row 131: The Mountain Steeped
$ws.Range("H131").Value = 20897.826
$ws.Range("I131").Value = 167020
$ws.Range("J131").Value = 1838.4131
$ws.Range("K131").Value = 501060
$ws.Range("L131").Value = 5515.2393
$ws.Range("M131").Value = -496020
$ws.Range("N131").Value = -15595.2393

$ws = $wb.Worksheets.Item("GSM")
# GSM row 70: Sky Is the Limit
$ws.Range("H70").Value = 4334.6206
$ws.Range("I70").Value = 4077.9546
$ws.Range("J70").Value = 5141.2856
$ws.Range("K70").Value = 4077.9546
$ws.Range("L70").Value = 5141.2856
$ws.Range("M70").Value = -3807.9546
$ws.Range("N70").Value = -5681.2856

# GSM row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 4334.6206
$ws.Range("I73").Value = 4077.9546
$ws.Range("J73").Value = 5141.2856
$ws.Range("K73").Value = 4077.9546
$ws.Range("L73").Value = 5141.2856
$ws.Range("M73").Value = -3141.9546
$ws.Range("N73").Value = -7013.2856

# GSM row 132: On Board for Lar
$ws.Range("H132").Value = 4328.4287
$ws.Range("I132").Value = 3550
$ws.Range("J132").Value = 5366.3335
$ws.Range("K132").Value = 10650
$ws.Range("L132").Value = 16099.0005
$ws.Range("M132").Value = -8120
$ws.Range("N132").Value = -21159.0005

# GSM row 138: Orders Anonymous
$ws.Range("H138").Value = 30943.334
$ws.Range("J138").Value = 30943.334
$ws.Range("L138").Value = 30943.334
$ws.Range("N138").Value = -41223.334

# GSM row 139: Ringing Gratitude
$ws.Range("H139").Value = 27192
$ws.Range("J139").Value = 27192
$ws.Range("L139").Value = 27192
$ws.Range("N139").Value = -37472

$ws = $wb.Worksheets.Item("LTW")
# LTW row 5: These Boots Are Made for Wailing
$ws.Range("H5").Value = 6966.6665
$ws.Range("I5").Value = 2950
$ws.Range("J5").Value = 15000
$ws.Range("K5").Value = 2950
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = -2837
$ws.Range("N5").Value = -15226

# LTW row 22: Skin off Their Backs
$ws.Range("H22").Value = 495.25
$ws.Range("I22").Value = 493.66666
$ws.Range("K22").Value = 493.66666
$ws.Range("M22").Value = -198.66666

# LTW row 27: Fire and Hide
$ws.Range("H27").Value = 495.25
$ws.Range("I27").Value = 493.66666
$ws.Range("K27").Value = 493.66666
$ws.Range("M27").Value = -386.66666

# LTW row 46: Supply Side Logic
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# LTW row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 248.92308
$ws.Range("I55").Value = 173.2
$ws.Range("J55").Value = 296.25
$ws.Range("K55").Value = 173.2
$ws.Range("L55").Value = 296.25
$ws.Range("M55").Value = -0.1999999999999886
$ws.Range("N55").Value = -642.25

# LTW row 61: Spelling Me Softly
$ws.Range("H61").Value = 14075.5
$ws.Range("I61").Value = 18100.666
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 18100.666
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -17898.666
$ws.Range("N61").Value = -2404

# LTW row 113: Peace in Rest
$ws.Range("H113").Value = 14075.5
$ws.Range("I113").Value = 18100.666
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 18100.666
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -15930.666
$ws.Range("N113").Value = -6340

# LTW row 134: Freezing Fingers
$ws.Range("H134").Value = 33476.332
$ws.Range("J134").Value = 33476.332
$ws.Range("L134").Value = 33476.332
$ws.Range("N134").Value = -43616.332

$ws = $wb.Worksheets.Item("WVR")
# WVR row 126: A Polished Purchase
$ws.Range("H126").Value = 744
$ws.Range("I126").Value = 687.1111
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 2061.3333
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = 408.6667000000002
$ws.Range("N126").Value = -7940

# WVR row 138: Halfgloves, Full Effort
$ws.Range("H138").Value = 39690
$ws.Range("J138").Value = 39690
$ws.Range("L138").Value = 39690
$ws.Range("N138").Value = -49970
